$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1591.7059
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 1622.4375
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 4867.3125
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -7083.3125
$ws.Range("H135").Value = 5080.613
$ws.Range("I135").Value = 777.3182
$ws.Range("J135").Value = 15599.777
$ws.Range("K135").Value = 6995.8638
$ws.Range("L135").Value = 140397.993
$ws.Range("M135").Value = -4460.8638
$ws.Range("N135").Value = -145467.993
$ws.Range("H138").Value = 3267.0657
$ws.Range("I138").Value = 2859.8333
$ws.Range("J138").Value = 3437.535
$ws.Range("K138").Value = 8579.499899999999
$ws.Range("L138").Value = 10312.605
$ws.Range("M138").Value = -3439.499899999999
$ws.Range("N138").Value = -20592.605
$ws.Range("H141").Value = 1786.1384
$ws.Range("I141").Value = 1207.5333
$ws.Range("J141").Value = 3088
$ws.Range("K141").Value = 3622.5999
$ws.Range("L141").Value = 9264
$ws.Range("M141").Value = 1557.4001
$ws.Range("N141").Value = -19624

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6175.78
$ws.Range("I32").Value = 4657.884
$ws.Range("J32").Value = 15500
$ws.Range("K32").Value = 4657.884
$ws.Range("L32").Value = 15500
$ws.Range("M32").Value = -4370.884
$ws.Range("N32").Value = -16074
$ws.Range("H45").Value = 15874010
$ws.Range("I45").Value = 33334218
$ws.Range("J45").Value = 1093.4546
$ws.Range("K45").Value = 33334218
$ws.Range("L45").Value = 1093.4546
$ws.Range("M45").Value = -33333841
$ws.Range("N45").Value = -1847.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5467.0293
$ws.Range("I20").Value = 5808.04
$ws.Range("J20").Value = 4519.778
$ws.Range("K20").Value = 5808.04
$ws.Range("L20").Value = 4519.778
$ws.Range("M20").Value = -5561.04
$ws.Range("N20").Value = -5013.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 500635
$ws.Range("I19").Value = 500635
$ws.Range("K19").Value = 500635
$ws.Range("M19").Value = -500465
$ws.Range("H24").Value = 500635
$ws.Range("I24").Value = 500635
$ws.Range("K24").Value = 500635
$ws.Range("M24").Value = -500465
$ws.Range("H31").Value = 3047.6667
$ws.Range("I31").Value = 2695.2917
$ws.Range("J31").Value = 5866.6665
$ws.Range("K31").Value = 2695.2917
$ws.Range("L31").Value = 5866.6665
$ws.Range("M31").Value = -2400.2917
$ws.Range("N31").Value = -6456.6665
$ws.Range("H34").Value = 3047.6667
$ws.Range("I34").Value = 2695.2917
$ws.Range("J34").Value = 5866.6665
$ws.Range("K34").Value = 2695.2917
$ws.Range("L34").Value = 5866.6665
$ws.Range("M34").Value = -2493.2917
$ws.Range("N34").Value = -6270.6665
$ws.Range("H62").Value = 22225028
$ws.Range("I62").Value = 2772.625
$ws.Range("J62").Value = 47621892
$ws.Range("K62").Value = 2772.625
$ws.Range("L62").Value = 47621892
$ws.Range("M62").Value = -2148.625
$ws.Range("N62").Value = -47623140
$ws.Range("H65").Value = 22225028
$ws.Range("I65").Value = 2772.625
$ws.Range("J65").Value = 47621892
$ws.Range("K65").Value = 13863.125
$ws.Range("L65").Value = 238109460
$ws.Range("M65").Value = -10743.125
$ws.Range("N65").Value = -238115700
$ws.Range("H99").Value = 2476.4119
$ws.Range("I99").Value = 2412.375
$ws.Range("J99").Value = 2533.3333
$ws.Range("K99").Value = 2412.375
$ws.Range("L99").Value = 2533.3333
$ws.Range("M99").Value = -914.375
$ws.Range("N99").Value = -5529.3333
$ws.Range("H122").Value = 998.3333
$ws.Range("I122").Value = 1021.6667
$ws.Range("J122").Value = 951.6667
$ws.Range("K122").Value = 3065.0001
$ws.Range("L122").Value = 2855.0001
$ws.Range("M122").Value = -615.0001000000002
$ws.Range("N122").Value = -7755.0001
$ws.Range("H126").Value = 2476.4119
$ws.Range("I126").Value = 2412.375
$ws.Range("J126").Value = 2533.3333
$ws.Range("K126").Value = 7237.125
$ws.Range("L126").Value = 7599.999899999999
$ws.Range("M126").Value = -4767.125
$ws.Range("N126").Value = -12539.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 384.65
$ws.Range("J107").Value = 522.2222
$ws.Range("L107").Value = 1566.6666
$ws.Range("N107").Value = -5406.6666
$ws.Range("H122").Value = 1077
$ws.Range("I122").Value = 849.75
$ws.Range("J122").Value = 1206.8572
$ws.Range("K122").Value = 7647.75
$ws.Range("L122").Value = 10861.7148
$ws.Range("M122").Value = -5197.75
$ws.Range("N122").Value = -15761.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32612832
$ws.Range("I70").Value = 42861236
$ws.Range("J70").Value = 4273.5454
$ws.Range("K70").Value = 42861236
$ws.Range("L70").Value = 4273.5454
$ws.Range("M70").Value = -42860966
$ws.Range("N70").Value = -4813.5454
$ws.Range("H73").Value = 32612832
$ws.Range("I73").Value = 42861236
$ws.Range("J73").Value = 4273.5454
$ws.Range("K73").Value = 42861236
$ws.Range("L73").Value = 4273.5454
$ws.Range("M73").Value = -42860300
$ws.Range("N73").Value = -6145.5454
$ws.Range("H97").Value = 1892.5927
$ws.Range("I97").Value = 1637.3914
$ws.Range("J97").Value = 3360
$ws.Range("K97").Value = 1637.3914
$ws.Range("L97").Value = 3360
$ws.Range("M97").Value = -1141.3914
$ws.Range("N97").Value = -4352
$ws.Range("H113").Value = 25001320
$ws.Range("J113").Value = 1466.6666
$ws.Range("L113").Value = 1466.6666
$ws.Range("N113").Value = -5806.6666
$ws.Range("H122").Value = 28573748
$ws.Range("I122").Value = 45456748
$ws.Range("J122").Value = 2515.3845
$ws.Range("K122").Value = 136370244
$ws.Range("L122").Value = 7546.1535
$ws.Range("M122").Value = -136367794
$ws.Range("N122").Value = -12446.1535
$ws.Range("H123").Value = 28879
$ws.Range("J123").Value = 28879
$ws.Range("L123").Value = 28879
$ws.Range("N123").Value = -33779
$ws.Range("H132").Value = 4374.325
$ws.Range("I132").Value = 4415.972
$ws.Range("K132").Value = 13247.916
$ws.Range("M132").Value = -10717.916

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 38463740
$ws.Range("I7").Value = 1957.7142
$ws.Range("J7").Value = 83335816
$ws.Range("K7").Value = 1957.7142
$ws.Range("L7").Value = 83335816
$ws.Range("M7").Value = -1845.7142
$ws.Range("N7").Value = -83336040
$ws.Range("H126").Value = 38463740
$ws.Range("I126").Value = 1957.7142
$ws.Range("J126").Value = 83335816
$ws.Range("K126").Value = 5873.142599999999
$ws.Range("L126").Value = 250007448
$ws.Range("M126").Value = -3403.142599999999
$ws.Range("N126").Value = -250012388
$ws.Range("H136").Value = 3482.6606
$ws.Range("I136").Value = 3586.7874
$ws.Range("J136").Value = 2938.889
$ws.Range("K136").Value = 10760.3622
$ws.Range("L136").Value = 8816.667000000001
$ws.Range("M136").Value = -8210.3622
$ws.Range("N136").Value = -13916.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2937.1304
$ws.Range("I136").Value = 3471.1667
$ws.Range("J136").Value = 2354.5454
$ws.Range("K136").Value = 10413.5001
$ws.Range("L136").Value = 7063.6362
$ws.Range("M136").Value = -7863.500100000001
$ws.Range("N136").Value = -12163.6362
